$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.449.76'
$ws.Range('E2').Value = '  +2.92%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.389.55'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.77'
$ws.Range('E5').Value = '  +2.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.80'
$ws.Range('E6').Value = '  +4.06%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.596'
$ws.Range('E8').Value = '  +1.14%  '
$ws.Range('E9').Value = '  +11.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.588'
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.81'
$ws.Range('E11').Value = '  +6.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000284'
$ws.Range('E12').Value = '  +5.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '688.15'
$ws.Range('E13').Value = '  -2.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.939.44'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.58'
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.512.58'
$ws.Range('E16').Value = '  +2.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.401.90'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.71'
$ws.Range('E19').Value = '  +2.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.39'
$ws.Range('E20').Value = '  +4.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.904'
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.40'
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.11'
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '104.52'
$ws.Range('E24').Value = '  +6.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.97'
$ws.Range('E25').Value = '  +3.19%  '
$ws.Range('E26').Value = '  +2.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.65'
$ws.Range('E27').Value = '  +3.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.48'
$ws.Range('E28').Value = '  +3.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.70'
$ws.Range('E29').Value = '  +2.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.05'
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.18'
$ws.Range('E31').Value = '  +2.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.66'
$ws.Range('E32').Value = '  +11.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '557.56'
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('E34').Value = '  +2.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.41'
$ws.Range('E35').Value = '  +2.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.728.72'
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.141'
$ws.Range('E38').Value = '  +8.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.00'
$ws.Range('E39').Value = '  +2.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0710'
$ws.Range('E40').Value = '  +7.20%  '
$ws.Range('E41').Value = '  +2.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.67'
$ws.Range('E42').Value = '  +2.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.340'
$ws.Range('E43').Value = '  +2.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0419'
$ws.Range('E44').Value = '  +3.94%  '
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.66'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('E47').Value = '  +1.75%  '
$ws.Range('E48').Value = '  +7.84%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.43'
$ws.Range('E50').Value = '  +2.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.65'
$ws.Range('E51').Value = '  -0.94%  '
